$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 64, pushing existing rows 64:126 down to 65:127.
$ws.Rows("64:64").Insert()

# The new blank row 64 should start life as a duplicate of the row that
# just landed at 65 (the original row 64's data), then get its own
# date / variety values.
$ws.Range("A65:T65").Copy()
$ws.Range("A64").PasteSpecial()

$ws.Range("D64").Value = 44874
$ws.Range("K64").Value = "Murcott"
